# Weekly data refresh: insert 6 new rows (one new sampling date) at the top
# of the Repollo / Vega Monumental Concepción data block (rows 397-402),
# pushing the existing rows down by 6. The oldest 6 rows that fall off the
# bottom of the historical window simply land at the new end of the sheet
# (rows 459-464) because Excel's row insert shifts everything below too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows right before row 397; everything from 397 down
# (including the previous 397-458 block) shifts down to 403-464.
$ws.Rows("397:402").Insert()

# New week's data: six Variedad x Calidad combinations, all dated 44951
# (2023-01-25), all in "Región Metropolitana".
$newRows = @(
    @{ Row = 397; Variedad = "Copenhague";    Calidad = "Primera"; J = 800;  K = 800; L = 900;  M = 850 },
    @{ Row = 398; Variedad = "Copenhague";    Calidad = "Segunda"; J = 400;  K = 700; L = 700;  M = 700 },
    @{ Row = 399; Variedad = "Crespo record"; Calidad = "Primera"; J = 1000; K = 900; L = 1000; M = 950 },
    @{ Row = 400; Variedad = "Crespo record"; Calidad = "Segunda"; J = 500;  K = 700; L = 700;  M = 700 },
    @{ Row = 401; Variedad = "Morada(o)";     Calidad = "Primera"; J = 500;  K = 900; L = 1000; M = 960 },
    @{ Row = 402; Variedad = "Morada(o)";     Calidad = "Segunda"; J = 300;  K = 800; L = 800;  M = 800 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = 11
    $ws.Cells.Item($row, 2).Value  = "Vega Monumental Concepción"
    $ws.Cells.Item($row, 3).Value  = "Bíobío"
    $ws.Cells.Item($row, 4).Value  = 44951
    $ws.Cells.Item($row, 5).Value  = 8
    $ws.Cells.Item($row, 6).Value  = 100112006
    $ws.Cells.Item($row, 7).Value  = "Repollo"
    $ws.Cells.Item($row, 8).Value  = $r.Variedad
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = "`$/unidad"
    $ws.Cells.Item($row, 15).Value = "Región Metropolitana"
    $ws.Cells.Item($row, 16).Value = $r.M
    $ws.Cells.Item($row, 17).Value = 1
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
